$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 10575.833
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 10575.833
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 31727.499
$ws.Cells.Item(70, 14).Value = -32267.499
$ws.Cells.Item(70, 13).ClearContents()

$ws.Cells.Item(73, 8).Value = 10575.833
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 10575.833
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 31727.499
$ws.Cells.Item(73, 14).Value = -33599.499
$ws.Cells.Item(73, 13).ClearContents()

$ws.Cells.Item(112, 8).Value = 2000

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 937.3333
$ws.Cells.Item(2, 9).Value = 949.5
$ws.Cells.Item(2, 11).Value = 949.5
$ws.Cells.Item(2, 13).Value = -836.5

$ws.Cells.Item(61, 8).Value = 12818.667
$ws.Cells.Item(61, 9).Value = 10637.333
$ws.Cells.Item(61, 11).Value = 10637.333
$ws.Cells.Item(61, 13).Value = -10425.333

$ws.Cells.Item(63, 8).Value = 4605.4
$ws.Cells.Item(63, 9).Value = 4605.4
$ws.Cells.Item(63, 11).Value = 4605.4
$ws.Cells.Item(63, 13).Value = -3919.4

$ws.Cells.Item(66, 8).Value = 4605.4
$ws.Cells.Item(66, 9).Value = 4605.4
$ws.Cells.Item(66, 11).Value = 23027
$ws.Cells.Item(66, 13).Value = -19595

$ws.Cells.Item(74, 8).Value = 5189.125
$ws.Cells.Item(74, 9).Value = 4859
$ws.Cells.Item(74, 10).Value = 7500
$ws.Cells.Item(74, 11).Value = 4859
$ws.Cells.Item(74, 12).Value = 7500
$ws.Cells.Item(74, 13).Value = -3985
$ws.Cells.Item(74, 14).Value = -9248

$ws.Cells.Item(77, 8).Value = 5189.125
$ws.Cells.Item(77, 9).Value = 4859
$ws.Cells.Item(77, 10).Value = 7500
$ws.Cells.Item(77, 11).Value = 24295
$ws.Cells.Item(77, 12).Value = 37500
$ws.Cells.Item(77, 13).Value = -19927
$ws.Cells.Item(77, 14).Value = -46236

$ws.Cells.Item(88, 8).Value = 2499.5
$ws.Cells.Item(88, 9).Value = 2524.25
$ws.Cells.Item(88, 10).Value = 2450
$ws.Cells.Item(88, 11).Value = 2524.25
$ws.Cells.Item(88, 12).Value = 2450
$ws.Cells.Item(88, 13).Value = -2118.25
$ws.Cells.Item(88, 14).Value = -3262

$ws.Cells.Item(91, 8).Value = 2499.5
$ws.Cells.Item(91, 9).Value = 2524.25
$ws.Cells.Item(91, 10).Value = 2450
$ws.Cells.Item(91, 11).Value = 2524.25
$ws.Cells.Item(91, 12).Value = 2450
$ws.Cells.Item(91, 13).Value = -1120.25
$ws.Cells.Item(91, 14).Value = -5258

$ws.Cells.Item(116, 8).Value = 937.3333
$ws.Cells.Item(116, 9).Value = 949.5
$ws.Cells.Item(116, 11).Value = 949.5
$ws.Cells.Item(116, 13).Value = 1344.5

$ws.Cells.Item(132, 8).Value = 7516
$ws.Cells.Item(132, 9).Value = 3880.1428
$ws.Cells.Item(132, 11).Value = 11640.4284
$ws.Cells.Item(132, 13).Value = -9110.428400000001

$ws.Cells.Item(136, 8).Value = 12818.667
$ws.Cells.Item(136, 9).Value = 10637.333
$ws.Cells.Item(136, 11).Value = 31911.999
$ws.Cells.Item(136, 13).Value = -29361.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 937.3333
$ws.Cells.Item(3, 9).Value = 949.5
$ws.Cells.Item(3, 11).Value = 949.5
$ws.Cells.Item(3, 13).Value = -835.5

$ws.Cells.Item(20, 8).Value = 4345.25
$ws.Cells.Item(20, 9).Value = 4345.25
$ws.Cells.Item(20, 11).Value = 4345.25
$ws.Cells.Item(20, 13).Value = -4098.25

$ws.Cells.Item(94, 8).Value = 4000
$ws.Cells.Item(94, 9).Value = 4000
$ws.Cells.Item(94, 11).Value = 4000
$ws.Cells.Item(94, 13).Value = -3549

$ws.Cells.Item(99, 8).Value = 1899.4286
$ws.Cells.Item(99, 9).Value = 1859.4
$ws.Cells.Item(99, 11).Value = 1859.4
$ws.Cells.Item(99, 13).Value = -361.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 10640
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 10640
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 10640
$ws.Cells.Item(58, 14).Value = -11046
$ws.Cells.Item(58, 13).ClearContents()

$ws.Cells.Item(105, 8).Value = 5803
$ws.Cells.Item(105, 9).Value = 5803
$ws.Cells.Item(105, 11).Value = 5803
$ws.Cells.Item(105, 13).Value = -4056

$ws.Cells.Item(136, 8).Value = 10640
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 10640
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 31920
$ws.Cells.Item(136, 14).Value = -37020
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 573.2222
$ws.Cells.Item(34, 10).Value = 1097.25
$ws.Cells.Item(34, 12).Value = 3291.75
$ws.Cells.Item(34, 14).Value = -3459.75

$ws.Cells.Item(108, 8).Value = 181.75
$ws.Cells.Item(108, 9).Value = 181.75
$ws.Cells.Item(108, 11).Value = 545.25
$ws.Cells.Item(108, 13).Value = 2334.75

$ws.Cells.Item(109, 8).Value = 2828.6667
$ws.Cells.Item(109, 9).Value = 2828.6667
$ws.Cells.Item(109, 11).Value = 8486.000100000001
$ws.Cells.Item(109, 13).Value = -7446.000100000001

$ws.Cells.Item(132, 8).Value = 1311.4546
$ws.Cells.Item(132, 9).Value = 1223.1666
$ws.Cells.Item(132, 11).Value = 11008.4994
$ws.Cells.Item(132, 13).Value = -8478.499400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 999
$ws.Cells.Item(70, 9).Value = 999
$ws.Cells.Item(70, 11).Value = 999
$ws.Cells.Item(70, 13).Value = -729

$ws.Cells.Item(73, 8).Value = 999
$ws.Cells.Item(73, 9).Value = 999
$ws.Cells.Item(73, 11).Value = 999
$ws.Cells.Item(73, 13).Value = -63

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1968.4286
$ws.Cells.Item(55, 9).Value = 1700
$ws.Cells.Item(55, 10).Value = 2169.75
$ws.Cells.Item(55, 11).Value = 1700
$ws.Cells.Item(55, 12).Value = 2169.75
$ws.Cells.Item(55, 13).Value = -1527
$ws.Cells.Item(55, 14).Value = -2515.75

$ws.Cells.Item(132, 8).Value = 11063.667
$ws.Cells.Item(132, 9).Value = 8095.5
$ws.Cells.Item(132, 11).Value = 24286.5
$ws.Cells.Item(132, 13).Value = -21756.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 2500
$ws.Cells.Item(14, 10).Value = 2500
$ws.Cells.Item(14, 12).Value = 2500
$ws.Cells.Item(14, 14).Value = -2836

$ws.Cells.Item(70, 8).Value = 25000
$ws.Cells.Item(70, 10).Value = 25000
$ws.Cells.Item(70, 12).Value = 25000
$ws.Cells.Item(70, 14).Value = -25630

$ws.Cells.Item(73, 8).Value = 25000
$ws.Cells.Item(73, 10).Value = 25000
$ws.Cells.Item(73, 12).Value = 25000
$ws.Cells.Item(73, 14).Value = -27184

$ws.Cells.Item(81, 8).Value = 8000
$ws.Cells.Item(81, 9).Value = 2000
$ws.Cells.Item(81, 11).Value = 4000
$ws.Cells.Item(81, 13).Value = -2939

$ws.Cells.Item(84, 8).Value = 8000
$ws.Cells.Item(84, 9).Value = 2000
$ws.Cells.Item(84, 11).Value = 20000
$ws.Cells.Item(84, 13).Value = -14696

$ws.Cells.Item(100, 8).Value = 1913
$ws.Cells.Item(100, 9).Value = 1650.3334
$ws.Cells.Item(100, 10).Value = 2701
$ws.Cells.Item(100, 11).Value = 3300.6668
$ws.Cells.Item(100, 12).Value = 5402
$ws.Cells.Item(100, 13).Value = -2759.6668
$ws.Cells.Item(100, 14).Value = -6484

$ws.Cells.Item(132, 8).Value = 12681.637
$ws.Cells.Item(132, 9).Value = 14999
$ws.Cells.Item(132, 11).Value = 44997
$ws.Cells.Item(132, 13).Value = -42467
